# Updated simulation files with Holden scheme
#
# - Adds a new "Holden" sampling scheme (2.5/5/10/15 deg) to the HKL sheet,
#   inserted where the old "HexGrid-90degTilt*" rows used to be; the
#   HexGrid rows move down to 4 brand-new rows at the bottom.
# - Drops the stray duplicate block of columns (X:AG) that repeated the
#   Pair/5A4F/MaxUnique headers a second time.
# - Re-orders the "[h, k, l]" HKL column headers in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the duplicated tail columns (X:AG) -------------------------------
$ws.Range("X1:AG19").EntireColumn.Delete()

# --- Re-order the HKL column headers across row 2 (C2:M2) ------------------
$ws.Range("C2").Value = "[3, 3, 1]"
$ws.Range("D2").Value = "[3, 1, 1]"
$ws.Range("E2").Value = "[1, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 2]"
$ws.Range("G2").Value = "[5, 1, 1]"
$ws.Range("H2").Value = "[4, 2, 2]"
$ws.Range("I2").Value = "[4, 2, 0]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "[3, 3, 3]"
$ws.Range("L2").Value = "[2, 0, 0]"
$ws.Range("M2").Value = "[2, 2, 0]"

# --- Relabel the old HexGrid-90degTilt* rows as the new Holden rows --------
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# --- Append the HexGrid-90degTilt* rows back at the bottom (20:23) ---------
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"

$ws.Range("C20:W23").Value = 1
